$wb = $excel.ActiveWorkbook
$wsSemicon = $wb.Worksheets.Item("semicon")
$wsAlloys  = $wb.Worksheets.Item("alloys")

# ---------------------------------------------------------------------------
# Sheet "semicon": insert one derived column (E = D/20) and add a small new
# row of scratch data (B8/C8).
# ---------------------------------------------------------------------------
$wsSemicon.Range("E1").EntireColumn.Insert()

for ($r = 2; $r -le 6; $r++) {
    $wsSemicon.Range("E$r").Formula = "=D$r/20"
}

$wsSemicon.Range("B8").Value = 2.86
$wsSemicon.Range("C8").Formula = "=B8*SQRT(2)/20"

# ---------------------------------------------------------------------------
# Sheet "alloys": insert 5 derived columns, one to the right of each of the
# existing computed columns (old E, F, G, H) plus one right after D.
# After each insert the previously-existing columns shift right, so the
# insert points below are always expressed in terms of the *final* layout:
#   D (intensity)            -> new E = D/20
#   F (old E, Kalph/Kbet)    -> new G = F*SQRT(2)/20
#   H (old F, k)             -> new I = H*SQRT(2)/20
#   J (old G, prop)          -> new K = J/10
#   L (old H, by mass)       -> new M = L/10
# ---------------------------------------------------------------------------
$wsAlloys.Range("E1").EntireColumn.Insert()
$wsAlloys.Range("G1").EntireColumn.Insert()
$wsAlloys.Range("I1").EntireColumn.Insert()
$wsAlloys.Range("K1").EntireColumn.Insert()
$wsAlloys.Range("M1").EntireColumn.Insert()

for ($r = 2; $r -le 28; $r++) {
    $wsAlloys.Range("E$r").Formula = "=D$r/20"
    $wsAlloys.Range("G$r").Formula = "=F$r*SQRT(2)/20"
    $wsAlloys.Range("I$r").Formula = "=H$r*SQRT(2)/20"
    $wsAlloys.Range("K$r").Formula = "=J$r/10"
    $wsAlloys.Range("M$r").Formula = "=L$r/10"
}

# Bit of scratch data the author pasted in below the table.
$wsAlloys.Range("J33").Value = 122.5

# ---------------------------------------------------------------------------
# Final selection / active sheet state.
# ---------------------------------------------------------------------------
$wsAlloys.Activate()
$wsAlloys.Range("J33").Select()

$wsSemicon.Activate()
$wsSemicon.Range("C9").Select()
